# Apply the "Added full e2e test of the add book flow" edit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New test-case row (row 3): delete-book / bc887d1000 ----------------
$ws.Range("A3").Value = "delete-book"
$ws.Range("B3").Value = "bc887d1000"

# --- New styled (but empty) cell at B18 ----------------------------------
# Font: Menlo, 12pt, color #CE9178 (a VS Code style orange string color).
# A temporary named style is used so both font properties (Name and
# Color) are resolved together before being applied to the cell; the
# temporary style is then removed again so it doesn't linger in the
# workbook's cell-style table, leaving only direct cell formatting behind.
$tmpStyleName = "MenloStyleTmp"
$tmpStyle = $wb.Styles.Add($tmpStyleName)
$tmpStyle.Font.Name = "Menlo"
$tmpStyle.Font.Size = 12
$tmpStyle.Font.Color = 7901646   # RGB(206, 145, 120) = 0xCE9178, BGR-packed for COM
$ws.Range("B18").Style = $tmpStyleName
$wb.Styles.Item($tmpStyleName).Delete()

# --- Update the sheet's active cell / selection to B4 --------------------
$ws.Range("B4").Select() | Out-Null
